# Cryptos list update - Tue Mar 19 17:31:16 UTC 2024
# D-column "Price" values look numeric (contain dots as thousands separators
# or plain decimals) but must stay as literal text, exactly as authored
# upstream. A leading apostrophe forces Excel to store them as text instead
# of auto-converting to a Number (which would also silently drop trailing
# zeros, e.g. "185.70" -> 185.7). The apostrophe itself is not stored as
# part of the cell's value/text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'65.720.57"
$ws.Range("E2").Value = "  -2.07%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.354.88"
$ws.Range("E3").Value = "  -3.85%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'185.70"
$ws.Range("E5").Value = "  -6.05%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'538.78"
$ws.Range("E6").Value = "  -0.84%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.608"
$ws.Range("E7").Value = "  +1.22%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "'3.343.39"
$ws.Range("E8").Value = "  -3.99%  "

# Row 9 - USDC (price unchanged)
$ws.Range("E9").Value = "  -0.08%  "

# Row 10 - was Cardano, now Avalanche
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").Value = "'61.46"
$ws.Range("E10").Value = "  -2.19%  "

# Row 11 - was Avalanche, now Cardano
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.628"
$ws.Range("E11").Value = "  -3.34%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "'0.136"
$ws.Range("E12").Value = "  -3.85%  "

# Row 13 - ShibaInu (price unchanged)
$ws.Range("E13").Value = "  -0.03%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'9.25"
$ws.Range("E14").Value = "  -4.73%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'3.867.18"
$ws.Range("E15").Value = "  -4.69%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "'3.340.68"
$ws.Range("E16").Value = "  -4.49%  "

# Row 17 - TRON
$ws.Range("D17").Value = "'0.119"
$ws.Range("E17").Value = "  -3.81%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'17.95"
$ws.Range("E18").Value = "  -1.93%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "'65.383.53"
$ws.Range("E19").Value = "  -1.95%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'11.31"
$ws.Range("E20").Value = "  -3.45%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "'0.975"
$ws.Range("E21").Value = "  -4.46%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'379.83"
$ws.Range("E22").Value = "  -2.34%  "

# Row 23 - PancakeSwap
$ws.Range("D23").Value = "'3.88"
$ws.Range("E23").Value = "  -2.06%  "

# Row 24 - RenderToken
$ws.Range("D24").Value = "'11.50"
$ws.Range("E24").Value = "  -2.09%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'82.22"
$ws.Range("E25").Value = "  +0.49%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "'3.91"
$ws.Range("E26").Value = "  +5.91%  "

# Row 27 - LEO
$ws.Range("D27").Value = "'6.06"
$ws.Range("E27").Value = "  -1.29%  "

# Row 28 - ImmutableX
$ws.Range("D28").Value = "'2.73"
$ws.Range("E28").Value = "  -1.69%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "'11.75"
$ws.Range("E29").Value = "  -2.35%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "'8.59"
$ws.Range("E30").Value = "  -1.27%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "'29.40"
$ws.Range("E31").Value = "  -4.03%  "

# Row 32 - Bittensor
$ws.Range("D32").Value = "'657.59"
$ws.Range("E32").Value = "  -2.77%  "

# Row 33 - NEARProtocol
$ws.Range("D33").Value = "'6.91"
$ws.Range("E33").Value = "  -1.42%  "

# Row 34 - Cosmos
$ws.Range("D34").Value = "'11.46"
$ws.Range("E34").Value = "  -1.30%  "

# Row 35 - was Hedera, now OKB
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'60.06"
$ws.Range("E35").Value = "  -4.93%  "

# Row 36 - was OKB, now Hedera
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = "  -1.99%  "

# Row 37 - Dai (price unchanged)
$ws.Range("E37").Value = "  -0.04%  "

# Row 38 - TheGraph (price unchanged)
$ws.Range("E38").Value = "  +0.56%  "

# Row 39 - InjectiveProtocol
$ws.Range("D39").Value = "'37.26"
$ws.Range("E39").Value = "  -3.23%  "

# Row 40 - PEPE
$ws.Range("D40").Value = "'0.0₃0735"
$ws.Range("E40").Value = "  +9.95%  "

# Row 41 - FirstDigitalUSD
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.35%  "

# Row 42 - Kaspa (price unchanged)
$ws.Range("E42").Value = "  -0.88%  "

# Row 43 - Maker
$ws.Range("D43").Value = "'2.929.91"
$ws.Range("E43").Value = "  -4.12%  "

# Row 44 - Fetch.AI
$ws.Range("D44").Value = "'2.54"
$ws.Range("E44").Value = "  +2.44%  "

# Row 45 - ThetaToken
$ws.Range("D45").Value = "'2.75"
$ws.Range("E45").Value = "  -7.78%  "

# Row 46 - VeChain
$ws.Range("D46").Value = "'0.0407"
$ws.Range("E46").Value = "  +3.32%  "

# Row 47 - was WEMIXToken, now Stacks
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "  +13.21%  "

# Row 48 - was Stacks, now WEMIXToken
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'2.68"
$ws.Range("E48").Value = "  -1.19%  "

# Row 49 - dogwifhat
$ws.Range("D49").Value = "'2.71"
$ws.Range("E49").Value = "  +1.85%  "

# Row 50 - Stellar (price unchanged)
$ws.Range("E50").Value = "  +1.49%  "

# Row 51 - ApeXProtocol
$ws.Range("D51").Value = "'3.01"
$ws.Range("E51").Value = "  +4.56%  "
